$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.982.31"
$ws.Range("E2").Value = "  +5.29%  "

$ws.Range("D3").Value = "1.880.12"
$ws.Range("E3").Value = "  +4.13%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "282.29"
$ws.Range("E5").Value = "  +2.15%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "0.5252"
$ws.Range("E7").Value = "  +3.78%  "

$ws.Range("D8").Value = "0.3545"
$ws.Range("E8").Value = "  +0.71%  "

$ws.Range("D9").Value = "45.24"
$ws.Range("E9").Value = "  +3.67%  "

$ws.Range("D10").Value = "0.07069"
$ws.Range("E10").Value = "  +6.56%  "

$ws.Range("D11").Value = "20.35"
$ws.Range("E11").Value = "  +1.66%  "

$ws.Range("D12").Value = "0.8198"
$ws.Range("E12").Value = "  -2.10%  "

$ws.Range("D13").Value = "0.07818"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").Value = "1.885.80"
$ws.Range("E14").Value = "  +4.48%  "

$ws.Range("D15").Value = "5.236"
$ws.Range("E15").Value = "  +3.16%  "

$ws.Range("D16").Value = "90.57"
$ws.Range("E16").Value = "  +3.50%  "

$ws.Range("E17").Value = "  +0.10%  "

$ws.Range("D18").Value = "14.60"
$ws.Range("E18").Value = "  +4.81%  "

$ws.Range("D19").Value = "0.000008159"
$ws.Range("E19").Value = "  +2.65%  "

$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "27.025.11"
$ws.Range("E21").Value = "  +5.21%  "

$ws.Range("D22").Value = "4.789"
$ws.Range("E22").Value = "  +1.48%  "

$ws.Range("D23").Value = "10.21"
$ws.Range("E23").Value = "  +1.90%  "

$ws.Range("D24").Value = "6.263"
$ws.Range("E24").Value = "  +3.62%  "

$ws.Range("D25").Value = "2.410"
$ws.Range("E25").Value = "  +14.32%  "

$ws.Range("D26").Value = "146.98"
$ws.Range("E26").Value = "  +3.09%  "

$ws.Range("D27").Value = "17.64"
$ws.Range("E27").Value = "  +4.35%  "

$ws.Range("D28").Value = "1.663"
$ws.Range("E28").Value = "  +0.54%  "

$ws.Range("D29").Value = "113.70"
$ws.Range("E29").Value = "  +4.88%  "

$ws.Range("D30").Value = "4.408"
$ws.Range("E30").Value = "  +2.17%  "

$ws.Range("D31").Value = "4.396"
$ws.Range("E31").Value = "  +4.28%  "

$ws.Range("D32").Value = "0.08884"
$ws.Range("E32").Value = "  +0.97%  "

$ws.Range("D33").Value = "0.04919"
$ws.Range("E33").Value = "  +2.76%  "

$ws.Range("D34").Value = "1.177"
$ws.Range("E34").Value = "  +4.66%  "

$ws.Range("D35").Value = "0.7476"
$ws.Range("E35").Value = "  +3.31%  "

$ws.Range("D36").Value = "2.903"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").Value = "3.297"
$ws.Range("E37").Value = "  +8.87%  "

$ws.Range("D38").Value = "2.400"
$ws.Range("E38").Value = "  +4.89%  "

$ws.Range("D39").Value = "0.5316"
$ws.Range("E39").Value = "  +3.30%  "

$ws.Range("D40").Value = "0.01895"
$ws.Range("E40").Value = "  +1.78%  "

$ws.Range("D41").Value = "0.9815"
$ws.Range("E41").Value = "  +2.09%  "

$ws.Range("D42").Value = "117.15"
$ws.Range("E42").Value = "  +2.24%  "

$ws.Range("D43").Value = "6.320"
$ws.Range("E43").Value = "  +2.37%  "

$ws.Range("D44").Value = "8.189"
$ws.Range("E44").Value = "  +1.99%  "

$ws.Range("D45").Value = "0.4635"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").Value = "0.1371"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("D48").Value = "9.487"
$ws.Range("E48").Value = "  +2.79%  "

$ws.Range("D49").Value = "36.87"
$ws.Range("E49").Value = "  +2.84%  "

$ws.Range("D50").Value = "1.529"
$ws.Range("E50").Value = "  +2.52%  "

$ws.Range("D51").Value = "0.05948"
$ws.Range("E51").Value = "  +2.37%  "

